# Rebuild the soft-skills list in column A of sheet 1 ("Лист1")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the previous (sparse, A1:A89) list entirely before rewriting it
$ws.Range("A1:A89").ClearContents()

$values = @(
  "soft_skills",
  "документация",
  "аналитические навыки",
  "коммуникация",
  "постановка задач разработчикам",
  "проактивность",
  "работа с большим объемом информации",
  "переговоры",
  "системное мышление",
  "сбор и анализ информации",
  "работа в команде",
  "сбор требований",
  "деловая коммуникация",
  "обучение",
  "организаторские навыки",
  "ответственность",
  "управление требованиями",
  "креативность",
  "грамотность",
  "удаленная работа",
  "управление персоналом",
  "проектный менеджмент",
  "ориентация на результат",
  "коммуникабельность",
  "внимание к деталям",
  "многозадачность",
  "документирование",
  "логика"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Narrow column A slightly (was 50.57 chars wide, now ~46.14)
$ws.Columns.Item(1).ColumnWidth = 45.3

# Move/restore the active selection to A30 (below the new data)
$ws.Range("A30").Select() | Out-Null

Write-Output "soft_skills sheet rebuilt"
